# Holding upload, does not create user verification email
#
# Inserts a new "Send Confirmation Email" column (Yes/No) after the existing
# "Email" column, updates the sample email domain from mycompany.com to
# myfirm.com, and removes the mailto hyperlinks that used to decorate the
# Email column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the mailto: hyperlinks that were attached to the Email column.
[void]$ws.Range("D2:D7").Hyperlinks.Delete()

# Insert a new column before the current First Name column (E), shifting
# First Name .. Grant Date one column to the right (F .. N).
$ws.Columns("E:E").Insert()

# The inserted column should carry the same (wide) column width that used
# to belong to the old column E ("First Name", now shifted to F).
$ws.Columns("E:E").ColumnWidth = 49.8366666666667

# New column header + values ("Yes"/"No" whether a confirmation email
# should be sent out for this holding).
$ws.Range("E1").Value = "Send Confirmation Email"
$ws.Range("E2:E5").Value = "No"
$ws.Range("E6:E7").Value = "Yes"

# Refresh the sample employee e-mail addresses to the new company domain.
$ws.Range("D2").Value = "emp1@myfirm.com"
$ws.Range("D3").Value = "emp2@myfirm.com"
$ws.Range("D4").Value = "emp3@myfirm.com"
$ws.Range("D5").Value = "emp4@myfirm.com"
$ws.Range("D6").Value = "emp5@myfirm.com"
$ws.Range("D7").Value = "emp6@myfirm.com"

# Match the selection left behind by the original edit.
[void]$ws.Range("D2:D7").Select()
